$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 and 38 swap content (THORChain <-> BinanceUSD), with updated price/volume values
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.01"
$ws.Range("E37").Value = "  +5.02%  "

$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.13%  "

# Row 45 and 46 swap content (Aave <-> TrustWalletToken), with updated price/volume values
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.18"
$ws.Range("E45").Value = "  +7.02%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "95.61"
$ws.Range("E46").Value = "  +6.46%  "

# Remaining price/volume updates across the list
$ws.Range("D2").Value = "37.527.83"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "2.077.33"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.29"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.48"
$ws.Range("E7").Value = "  +6.56%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.17"
$ws.Range("E10").Value = "  +3.34%  "
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "2.383.30"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.52"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.08"
$ws.Range("E15").Value = "  +4.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.782"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.21"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "2.087.77"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "37.678.82"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.21"
$ws.Range("E20").Value = "  +16.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.57"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "0.0₃0817"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.28"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  +5.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.09"
$ws.Range("E27").Value = "  +2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.49"
$ws.Range("E28").Value = "  +7.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.24"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0627"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("E35").Value = "  +7.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  +8.32%  "
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.54"
$ws.Range("E41").Value = "  +18.91%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0953"
$ws.Range("E43").Value = "  +3.12%  "
$ws.Range("D44").Value = "1.470.09"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.87"
$ws.Range("E48").Value = "  +3.36%  "
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.24"
$ws.Range("E50").Value = "  +5.62%  "
$ws.Range("E51").Value = "  +1.96%  "
